$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 175; existing rows 175-269 shift down to 176-270.
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new record (Brasil, 2022-11-08).
$ws.Cells.Item(175, 1).Value = 4
$ws.Cells.Item(175, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(175, 3).Value = "Los Lagos"
$ws.Cells.Item(175, 4).Value = 44873
$ws.Cells.Item(175, 5).Value = 10
$ws.Cells.Item(175, 6).Value = "Fruta"
$ws.Cells.Item(175, 7).Value = 100108
$ws.Cells.Item(175, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(175, 9).Value = 100108002
$ws.Cells.Item(175, 10).Value = "Mango"
$ws.Cells.Item(175, 11).Value = "Sin especificar"
$ws.Cells.Item(175, 12).Value = "Primera"
$ws.Cells.Item(175, 13).Value = 200
$ws.Cells.Item(175, 14).Value = 8000
$ws.Cells.Item(175, 15).Value = 9000
$ws.Cells.Item(175, 16).Value = 8500
$ws.Cells.Item(175, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(175, 18).Value = "Brasil"
$ws.Cells.Item(175, 19).Value = 2125
$ws.Cells.Item(175, 20).Value = 4
